# de_dg_gwq_lower.xlsx parameter update
# "add parameter comparison code, change base and v4 parameterizations"
#
# For this base parameter file, the substantive change is clearing the
# override values that had been entered for the "investor growers (white
# area)" row (row 18) in the C:F (DE, DG, GWQ, Lower-ish scenario) columns.
# The cells keep their existing number-format/style, only the values are
# removed, so ClearContents (not ClearFormats/Clear) is used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

[void]$ws.Range("C18:F18").ClearContents()

# Reflect the final cell selection left active on the sheet after the edit.
[void]$ws.Range("I21").Select()
